$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the empty paragraph after "Explain the action taken..." with the
#    full two-phase-commit restart-recovery answer (several new paragraphs).
# ---------------------------------------------------------------------------
$para4Xml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">There are two parts to answering this question. The first section will answer what happens when the one of the subordinate nodes fails. The second part will answer what is done if one of the coordinator </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>node</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> fails.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">If a subordinate node fails then there are several steps we may need to take for T depending on what is in the log. If we have an end record for T then we don’t need to do anything. However, if we have </w:t></w:r><w:r><w:t>a log record for T but we don’t have an end record than we will need to undo or redo T. If we have a prepare log for T but not a commit or abort we can say that this node is a subordinate. We will then contact the coordinator to enquire about the status of T and will then accordingly redo or undo T. Additionally, if we have no prepare log for T then we need to abort T.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">If we are the coordinator and we have a log for T and we have an end log for T we don’t need to worry about T because T was completed. However, if we do not have an end log for T then we need to keep sending commit/abort messages to subordinates until the acknowledge the request. Again if we have no prepare log for T we need to abort T along with sending out the abort message to all subordinates for T. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>
'@

$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML($para4Xml)

Write-Output "step1 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 2) "Again, assume ..." gets a mid-run bookmark (_GoBack moves here) and a
#    spell-check proofErr wrap around "muli".
# ---------------------------------------------------------------------------
$para14Xml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Ag</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">ain, assume that T1 reads and T2 and T3 read and update only data item D provide a schedule where the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>muli</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-Version concurrency control method will restart T2.</w:t></w:r></w:p>
'@
$p14 = $d.Paragraphs.Item(14)
$p14.Range.InsertXML($para14Xml)
Write-Output "step2 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 3) "Eventual consistency and vector clocks:" gets a lastRenderedPageBreak
#    marker before its run text.
# ---------------------------------------------------------------------------
$para19Xml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Eventual consistency and vector clocks:</w:t></w:r></w:p>
'@
$p19 = $d.Paragraphs.Item(19)
$p19.Range.InsertXML($para19Xml)
Write-Output "step3 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 4) "Assume the Following graph depict ..." gets grammar/spelling proofErr
#    wraps around "graph depict" and "preform".
# ---------------------------------------------------------------------------
$para22Xml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Assume the Following </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>graph depict</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> a part of the Web, where nodes represent pages and edges show hyper-links. Find out the pages whose PageRank values are greater than-zero and their relative PageRank values in the graph. You do not need to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>preform</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the fix point computation to determine the PageRank values. Instead, you should guess the PageRank values based on your understanding of the PageRank algorithm and explain why you think they are correct. If it is not possible to make any educated guess for some page(s), you should explain why.</w:t></w:r></w:p>
'@
$p22 = $d.Paragraphs.Item(22)
$p22.Range.InsertXML($para22Xml)
Write-Output "step4 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 5) "MongDB " gets a spell-check proofErr wrap around "MongDB".
# ---------------------------------------------------------------------------
$para24Xml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>MongDB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$p24 = $d.Paragraphs.Item(24)
$p24.Range.InsertXML($para24Xml)
Write-Output "step5 done, count=$($d.Paragraphs.Count)"

# ---------------------------------------------------------------------------
# 6) The trailing "_GoBack" bookmark on the last paragraph is removed (it
#    moved to paragraph 14 in step 2 above).
# ---------------------------------------------------------------------------
$para28Xml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Write an aggregate query to show total count of restaurants in each borough.</w:t></w:r></w:p>
'@
$p28 = $d.Paragraphs.Item(28)
$p28.Range.InsertXML($para28Xml)
Write-Output "step6 done, count=$($d.Paragraphs.Count)"
